# Update Window-to-Wall Ratio (WWR) values for SG Office, Retail, Restaurant
# archetypes on the ARCHITECTURE sheet: raise wwr_north/south/east/west
# from 0.35 to 0.59 (Sinberbest benchmark office), for rows:
#   row 5 -> OFFICE
#   row 6 -> RETAIL
#   row 8 -> RESTAURANT
# Columns G:J correspond to wwr_north, wwr_south, wwr_east, wwr_west.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

foreach ($row in 5, 6, 8) {
    $ws.Range("G" + $row + ":J" + $row).Value = 0.59
}

# Match the author's updated active-cell selection on the sheet view.
$ws.Activate()
$ws.Range("G8").Select()
